# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.841.01'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.293.04'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '112.55'
$ws.Range('E5').Value = '  +16.80%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.90'
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.05'
$ws.Range('E10').Value = '  +6.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +1.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.12'
$ws.Range('E12').Value = '  +15.87%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.81'
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.638.29'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.286.56'
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.719.11'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.73'
$ws.Range('E20').Value = '  +8.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.29'
$ws.Range('E21').Value = '  +0.43%  '
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.81'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.80'
$ws.Range('E24').Value = '  +7.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.86'
$ws.Range('E25').Value = '  +7.99%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.71'
$ws.Range('E27').Value = '  +4.42%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '41.90'
$ws.Range('E28').Value = '  +9.01%  '
$ws.Range('E29').Value = '  -2.06%  '
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.90'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0931'
$ws.Range('E32').Value = '  +3.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.56'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.69'
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('E37').Value = '  +3.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.107'
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.84'
$ws.Range('E39').Value = '  +14.17%  '
$ws.Range('E40').Value = '  +14.63%  '
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('E42').Value = '  +13.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.40'
$ws.Range('E43').Value = '  +3.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.36'
$ws.Range('E44').Value = '  +23.57%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('E46').Value = '  +3.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.78'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.98'
$ws.Range('E48').Value = '  +5.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0999'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.23'
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.467'
$ws.Range('E51').Value = '  +6.11%  '
